$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 — this shifts the existing rows 23..149
# down to 24..150, preserving their data/formatting exactly (matches the
# diff, where every row N (24..150) ends up holding what used to be in
# row N-1, and the sheet's used range grows from A1:R149 to A1:R150).
$ws.Rows("23:23").Insert()

# Populate the newly-inserted row 23 with the new data record.
$ws.Cells.Item(23, 1).Value = 1
$ws.Cells.Item(23, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(23, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(23, 4).Value = 45243
$ws.Cells.Item(23, 5).Value = 15
$ws.Cells.Item(23, 6).Value = 100112038
$ws.Cells.Item(23, 7).Value = "Cebollín baby"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 300
$ws.Cells.Item(23, 11).Value = 2500
$ws.Cells.Item(23, 12).Value = 3000
$ws.Cells.Item(23, 13).Value = 2750
$ws.Cells.Item(23, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(23, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 16).Value = 1375
$ws.Cells.Item(23, 17).Value = 2
$ws.Cells.Item(23, 18).Value = "Hortaliza"
